$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.552.37"
$ws.Range("E2").Value = "'  +2.48%  "

$ws.Range("D3").Value = "'1.671.60"
$ws.Range("E3").Value = "'  +2.08%  "

$ws.Range("E4").Value = "'  +0.24%  "

$ws.Range("D5").Value = "'239.25"
$ws.Range("E5").Value = "'  +1.28%  "

$ws.Range("E6").Value = "'  -0.01%  "

$ws.Range("D7").Value = "'0.4772"
$ws.Range("E7").Value = "'  +1.05%  "

$ws.Range("D8").Value = "'0.2619"
$ws.Range("E8").Value = "'  +2.66%  "

$ws.Range("D9").Value = "'0.06180"
$ws.Range("E9").Value = "'  +2.94%  "

$ws.Range("D10").Value = "'1.673.86"
$ws.Range("E10").Value = "'  +2.27%  "

$ws.Range("D11").Value = "'0.06983"
$ws.Range("E11").Value = "'  -0.47%  "

$ws.Range("D12").Value = "'14.87"
$ws.Range("E12").Value = "'  +0.51%  "

$ws.Range("D13").Value = "'0.5901"
$ws.Range("E13").Value = "'  -4.11%  "

$ws.Range("D14").Value = "'4.381"
$ws.Range("E14").Value = "'  +0.81%  "

$ws.Range("D15").Value = "'75.42"
$ws.Range("E15").Value = "'  +3.75%  "

$ws.Range("E16").Value = "'  +0.02%  "

$ws.Range("D17").Value = "'1.0000"
$ws.Range("E17").Value = "'  +0.15%  "

$ws.Range("D18").Value = "'25.539.95"
$ws.Range("E18").Value = "'  +2.39%  "

$ws.Range("D19").Value = "'0.000006753"
$ws.Range("E19").Value = "'  +2.68%  "

$ws.Range("D20").Value = "'11.46"
$ws.Range("E20").Value = "'  +3.14%  "

$ws.Range("D21").Value = "'1.887.66"
$ws.Range("E21").Value = "'  +2.19%  "

$ws.Range("D22").Value = "'4.457"
$ws.Range("E22").Value = "'  +1.44%  "

$ws.Range("D23").Value = "'8.809"
$ws.Range("E23").Value = "'  +2.53%  "

$ws.Range("D24").Value = "'5.274"
$ws.Range("E24").Value = "'  +0.20%  "

$ws.Range("D25").Value = "'136.74"
$ws.Range("E25").Value = "'  +2.54%  "

$ws.Range("E26").Value = "'  +1.83%  "

$ws.Range("D27").Value = "'1.390"
$ws.Range("E27").Value = "'  +1.73%  "

$ws.Range("D28").Value = "'1.729"
$ws.Range("E28").Value = "'  +4.34%  "

$ws.Range("D29").Value = "'104.67"
$ws.Range("E29").Value = "'  +2.18%  "

$ws.Range("D30").Value = "'3.987"
$ws.Range("E30").Value = "'  +6.53%  "

$ws.Range("D31").Value = "'0.07867"
$ws.Range("E31").Value = "'  +1.97%  "

$ws.Range("D32").Value = "'3.637"
$ws.Range("E32").Value = "'  +2.29%  "

$ws.Range("D33").Value = "'0.9991"
$ws.Range("E33").Value = "'  +0.04%  "

$ws.Range("D34").Value = "'0.04294"
$ws.Range("E34").Value = "'  -0.17%  "

$ws.Range("E35").Value = "'  +0.78%  "

$ws.Range("D36").Value = "'0.9569"
$ws.Range("E36").Value = "'  +4.03%  "

$ws.Range("D37").Value = "'0.6084"
$ws.Range("E37").Value = "'  +4.72%  "

$ws.Range("D38").Value = "'2.591"
$ws.Range("E38").Value = "'  +1.03%  "

$ws.Range("D39").Value = "'0.8934"
$ws.Range("E39").Value = "'  +7.78%  "

$ws.Range("E40").Value = "'  +0.18%  "

$ws.Range("E41").Value = "'  +3.57%  "

$ws.Range("D42").Value = "'0.01481"
$ws.Range("E42").Value = "'  -4.46%  "

$ws.Range("D43").Value = "'96.28"
$ws.Range("E43").Value = "'  -1.09%  "

$ws.Range("D44").Value = "'0.3762"
$ws.Range("E44").Value = "'  +1.38%  "

$ws.Range("D45").Value = "'4.927"
$ws.Range("E45").Value = "'  +4.02%  "

$ws.Range("D46").Value = "'0.1120"
$ws.Range("E46").Value = "'  +1.63%  "

$ws.Range("D47").Value = "'6.229"
$ws.Range("E47").Value = "'  +2.52%  "

$ws.Range("D48").Value = "'0.05267"
$ws.Range("E48").Value = "'  +1.02%  "

$ws.Range("D49").Value = "'29.97"
$ws.Range("E49").Value = "'  +1.44%  "

$ws.Range("D50").Value = "'7.423"
$ws.Range("E50").Value = "'  +3.31%  "

$ws.Range("E51").Value = "'  +0.21%  "
